$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.685.36"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").Value = "2.099.17"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0783"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.50%  "

$ws.Range("E11").Value = "  +1.27%  "

$ws.Range("D12").Value = "2.408.41"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.791"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").Value = "2.079.52"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").Value = "37.615.53"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("E28").Value = "  +4.10%  "

$ws.Range("E29").Value = "  -4.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0622"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.50%  "

$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0962"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.89%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.484.08"
$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("E44").Value = "  +0.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.83%  "

$ws.Range("E47").Value = "  +0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.01%  "

$ws.Range("D51").Value = "2.294.47"
$ws.Range("E51").Value = "  +1.07%  "
